$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; unprotect for editing, then restore protection afterwards.
$ws.Unprotect()

# Update the "as of" date in the confidential disclosure note (2021-05-05 -> 2021-05-06).
$ws.Range("A59").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-05-06 for illustrative purposes only and are subject to change."

# Update Weight (column D) and Percent Change (column E) values for holdings rows 2-56.
$ws.Range("D2").Value = 0.01525619417883049
$ws.Range("E2").Value = 0.001542614731970815
$ws.Range("D3").Value = 0.05090627932919957
$ws.Range("E3").Value = 0.01095537739945063
$ws.Range("D4").Value = 0.01456177284596408
$ws.Range("E4").Value = 0.005066592556170813
$ws.Range("D5").Value = 0.009574403780714266
$ws.Range("E5").Value = 0.01261542463259202
$ws.Range("D6").Value = 0.01522970238010053
$ws.Range("E6").Value = 0.0115713392174579
$ws.Range("D7").Value = 0.01952679042885419
$ws.Range("E7").Value = 0.01258804136070713
$ws.Range("D8").Value = 0.004460864021865693
$ws.Range("E8").Value = 0.0115773533290997
$ws.Range("D9").Value = 0.006832829479244817
$ws.Range("E9").Value = 0.009512875184516822
$ws.Range("D10").Value = 0.01446330802882791
$ws.Range("E10").Value = 0.008562075044069362
$ws.Range("D11").Value = 0.008440641959626955
$ws.Range("E11").Value = 0.004130707383639631
$ws.Range("D12").Value = 0.01552151685870768
$ws.Range("E12").Value = 0.00890493381468116
$ws.Range("D13").Value = 0.002889505003660743
$ws.Range("E13").Value = -0.02149321266968318
$ws.Range("D14").Value = 0.006002836132901327
$ws.Range("E14").Value = -0.0233884768967485
$ws.Range("D15").Value = 0.01422052361250007
$ws.Range("E15").Value = 0.02012442864398167
$ws.Range("D16").Value = 0.01032573111602295
$ws.Range("E16").Value = 0.01451634784968125
$ws.Range("D17").Value = 0.02086299192930808
$ws.Range("E17").Value = -0.003558718861210064
$ws.Range("D18").Value = 0.008456705141935607
$ws.Range("E18").Value = -0.003401360544217913
$ws.Range("D19").Value = 0.01664295112127905
$ws.Range("E19").Value = 0.004010295085892235
$ws.Range("D20").Value = 0.01200442705787175
$ws.Range("E20").Value = 0.005601369223587893
$ws.Range("D21").Value = 0.00737854185496312
$ws.Range("E21").Value = 0.01186387761473595
$ws.Range("D22").Value = 0.01447628932150758
$ws.Range("E22").Value = 0.007947976878612595
$ws.Range("D23").Value = 0.01980758481727284
$ws.Range("E23").Value = 0.005940762114911236
$ws.Range("D24").Value = 0.01011469950193314
$ws.Range("E24").Value = 0.01865405212424132
$ws.Range("D25").Value = 0.02003228258840431
$ws.Range("E25").Value = 0.01268686868686864
$ws.Range("D26").Value = 0.01402390528021047
$ws.Range("E26").Value = 0.00121200820436318
$ws.Range("D27").Value = 0.02092260625900006
$ws.Range("E27").Value = 0.01054009819967261
$ws.Range("D28").Value = 0.05582889758204242
$ws.Range("E28").Value = 0.01280249804839984
$ws.Range("D29").Value = 0.02134569678372268
$ws.Range("E29").Value = 0.005725611098877037
$ws.Range("D30").Value = 0.02966707895376738
$ws.Range("E30").Value = 0.008226652675760615
$ws.Range("D31").Value = 0.01550922043039001
$ws.Range("E31").Value = 0.01002599331600451
$ws.Range("D32").Value = 0.01342408861989442
$ws.Range("E32").Value = 0.01110327811068035
$ws.Range("D33").Value = 0.01908530195695247
$ws.Range("E33").Value = -0.004110393423370606
$ws.Range("D34").Value = 0.04323554943325002
$ws.Range("E34").Value = 0.009754748851937833
$ws.Range("D35").Value = 0.01080342400851567
$ws.Range("E35").Value = -0.002766251728907321
$ws.Range("D36").Value = 0.01009191219679761
$ws.Range("E36").Value = 0.002591121091725546
$ws.Range("D37").Value = 0.01055313721440414
$ws.Range("E37").Value = -0.005973451327433654
$ws.Range("D38").Value = 0.007272325618457082
$ws.Range("E38").Value = 0.006677796327211771
$ws.Range("D39").Value = 0.01198285383046885
$ws.Range("E39").Value = 0.01497946363856006
$ws.Range("D40").Value = 0.0168718825993756
$ws.Range("E40").Value = 0.02562096616467802
$ws.Range("D41").Value = 0.01720952072983998
$ws.Range("E41").Value = 0.006598845202089754
$ws.Range("D42").Value = 0.0322085483407407
$ws.Range("E42").Value = -0.004900255161215394
$ws.Range("D43").Value = 0.0113245460180145
$ws.Range("E43").Value = 0.00760877286078987
$ws.Range("D44").Value = 0.02211959351277898
$ws.Range("E44").Value = 0.00920553204484964
$ws.Range("D45").Value = 0.01243486430938669
$ws.Range("E45").Value = 0.01070229592411498
$ws.Range("D46").Value = 0.008588105708844159
$ws.Range("E46").Value = 0.01445209278047832
$ws.Range("D47").Value = 0.01332945281714579
$ws.Range("E47").Value = -0.003405078190684363
$ws.Range("D48").Value = 0.01034170090773679
$ws.Range("E48").Value = 0.00812743823146933
$ws.Range("D49").Value = 0.01589526601916826
$ws.Range("E49").Value = 0.01890694586432606
$ws.Range("D50").Value = 0.008469530783623908
$ws.Range("E50").Value = 0.01430892276930784
$ws.Range("D51").Value = 0.01203175937195507
$ws.Range("E51").Value = -0.002064693737095502
$ws.Range("D52").Value = 0.008552928584796147
$ws.Range("E52").Value = -0.03587300316290998
$ws.Range("D53").Value = 0.01002588504626147
$ws.Range("E53").Value = 0.01201938751113918
$ws.Range("D54").Value = 0.1351691887114325
$ws.Range("E54").Value = 0.0000985707244947065
$ws.Range("D55").Value = 0.04369185587952949
$ws.Range("E55").Value = 0.006177325581395277
$ws.Range("E56").Value = 0.006234162145849265

# Restore sheet protection (content protected; matches original workbook state).
$ws.Protect()
